$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Configuration")

# --- Remove the old "Mandetory_Files" (row10) / "Optional_Files" (row11) rows,
# and the old "MandetoryFiles" (ends up at row12 after the first delete) row,
# shifting everything below them up. This leaves (in order):
#   row9  = Name_Of_File / Invalid
#   row10 = Only_Generate_Mandetory_Files / Yes        (style carried from old row12)
#   row11 = InValidNameFilesInclude / Mandetory        (style carried from old row13)
#   row12..14 = blank filler rows
$ws.Range("A10:A11").EntireRow.Delete()
$ws.Range("A12").EntireRow.Delete()

# --- Update Name_Of_File value
$ws.Range("B9").Value = "Valid"

# --- Swap the text of the two consolidated rows (styles already match target
# because row10/row11 inherited s="7"/s="5" from the old row12/row13) and set
# the new values for them.
$ws.Range("A10").Value = "InValidNameFilesInclude"
$ws.Range("B10").Value = "Optional"
$ws.Range("A11").Value = "Only_Generate_Mandetory_Files"
$ws.Range("B11").Value = "Yes"

# --- New explanatory notes in column C
$ws.Range("C10").Value = "This is Applicable when Name of File is Invalid."
$ws.Range("C11").Value = "This is Applicable when Name of File is Valid."

# --- Insert one more blank filler row (with the same style as the existing
# blank rows) so the sheet ends at row 15 instead of row 14.
$ws.Range("A15").EntireRow.Insert()
$ws.Range("A15:B15").Style = $ws.Range("A14:B14").Style

# --- Column C width
$ws.Columns.Item(3).ColumnWidth = 63.140625

# --- Data validations: drop the old three and recreate them against the new
# ranges.
$ws.Range("B3:B12").Validation.Delete()

$ws.Range("B9:B10").Validation.Add(3, 1, 1, "Valid,Invalid")
$ws.Range("B9:B10").Validation.IgnoreBlank = $true
$ws.Range("B9:B10").Validation.InCellDropdown = $true
$ws.Range("B9:B10").Validation.ShowInput = $true
$ws.Range("B9:B10").Validation.ShowError = $true

$ws.Range("B3:B8,B11").Validation.Add(3, 1, 1, "Yes,No")
$ws.Range("B3:B8,B11").Validation.IgnoreBlank = $true
$ws.Range("B3:B8,B11").Validation.InCellDropdown = $true
$ws.Range("B3:B8,B11").Validation.ShowInput = $true
$ws.Range("B3:B8,B11").Validation.ShowError = $true

$ws.Range("B10").Validation.Add(3, 1, 1, "Mandetory,Optional,Both")
$ws.Range("B10").Validation.IgnoreBlank = $true
$ws.Range("B10").Validation.InCellDropdown = $true
$ws.Range("B10").Validation.ShowInput = $true
$ws.Range("B10").Validation.ShowError = $true

# --- Selection, matching the saved cursor position from the edit.
$ws.Range("A8").Select()

Write-Output "done"
